# Generate Report for Handoff
# Updates localization status from "In Translation" to "Ready for handoff"
# and refreshes the related handoff timestamps, then resizes the affected
# Status columns to fit the new (longer) text.

$wb = $excel.ActiveWorkbook

$oldStatus = "In Translation"
$newStatus = "Ready for handoff"

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
# E2/F2 hold the per-language status ("In Translation" -> "Ready for handoff")
$ws_overview.Range("E2").Value = $newStatus
$ws_overview.Range("F2").Value = $newStatus
# G2 holds "Latest HO Xliff Generate Date" - bump it forward
$ws_overview.Range("G2").Value = "2016-08-25 17:01:32"

# --- zh-cn sheet ---
# C2 holds the Status column
$ws_zhcn.Range("C2").Value = $newStatus
# H2 holds "Latest Handoff Datetime"
$ws_zhcn.Range("H2").Value = "2016-08-25 17:01:28"

# --- de-de sheet ---
# C2 holds the Status column
$ws_dede.Range("C2").Value = $newStatus
# H2 holds "Latest Handoff Datetime"
$ws_dede.Range("H2").Value = "2016-08-25 17:01:32"

# --- Resize the columns that display the status text so the longer
#     "Ready for handoff" string fits (matches the column autosize that
#     Excel performs after the content changes). ---
$ws_overview.Columns.Item(5).ColumnWidth = 16.3
$ws_overview.Columns.Item(6).ColumnWidth = 16.3

$ws_zhcn.Columns.Item(3).ColumnWidth = 16.3

$ws_dede.Columns.Item(3).ColumnWidth = 16.3
